$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("1000 Bs = 6.56 = 25973.77 pesos", "1000 Bs = 6.5 = 25793.24 pesos")
$text = $text.Replace("25973.77 pesos = 6.53 = 960.54 Bs", "25793.24 pesos = 6.49 = 971.32 Bs")
$cell.Value2 = $text

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 153.8
$ws2.Range("O10").Value = 3967
$ws2.Range("N12").Value = 3974.99
$ws2.Range("O12").Value = 149.69
